$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-InlineLogo($story, $newName) {
    $shape = $story.Range.InlineShapes.Item(1)
    # Selecting the shape first routes the subsequent property write through
    # $word.Selection.InlineShapes, which (unlike the Headers/Footers collection
    # range directly) reliably resolves the underlying picture for both header-
    # and footer-hosted inline pictures in this host.
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Pearson logo (footers) : image2.png -> image1.png
Rename-InlineLogo $sec.Footers.Item(1) "image1.png"
Rename-InlineLogo $sec.Footers.Item(2) "image1.png"

# BTEC logo (headers) : image1.jpg -> image2.jpg
Rename-InlineLogo $sec.Headers.Item(1) "image2.jpg"
Rename-InlineLogo $sec.Headers.Item(2) "image2.jpg"
